$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A labels for existing rows 2-10 (replace old "Conjoint: ..." labels) ---
$ws.Cells.Item(2, 1).Value  = "econ1"
$ws.Cells.Item(3, 1).Value  = "econ2"
$ws.Cells.Item(4, 1).Value  = "econ3"
$ws.Cells.Item(5, 1).Value  = "econ4"
$ws.Cells.Item(6, 1).Value  = "soc1"
$ws.Cells.Item(7, 1).Value  = "soc2"
$ws.Cells.Item(8, 1).Value  = "climate1"
$ws.Cells.Item(9, 1).Value  = "climate2"
$ws.Cells.Item(10, 1).Value = "climate3"

# --- New rows 11-16 with new labels ---
$ws.Cells.Item(11, 1).Value = "tax1: National redistribution scheme"
$ws.Cells.Item(12, 1).Value = "tax2: Wealth tax"
$ws.Cells.Item(13, 1).Value = "foreign1: Global climate scheme"
$ws.Cells.Item(14, 1).Value = "foreign2: Global tax on millionaires"
$ws.Cells.Item(15, 1).Value = "foreign3: Global democratic assembly on climate change"
$ws.Cells.Item(16, 1).Value = "foreign4: Doubling foreign aid"

# --- Updated numeric data for columns B:F, rows 2-16 ---
$data = @{
  2  = @(0.616858237547893, 0.676258992805755, 0.960199004975124, 0.953307392996109, 0.655973451327434)
  3  = @(0.845559845559846, 0.844594594594595, 0.869158878504673, 0.92057761732852,  0.840262582056893)
  4  = @(0.539748953974895, 0.815181518151815, 0.802690582959641, 0.61003861003861,  0.803474484256243)
  5  = @(0.895752895752896, 0.807817589576547, 0.934579439252336, 0.881294964028777, 0.88628762541806)
  6  = @(0.686507936507937, 0.762345679012346, 0.762931034482759, 0.807547169811321, 0.574786324786325)
  7  = @(0.697247706422018, 0.545454545454545, 0.684210526315789, 0.428044280442804, 0.598941798941799)
  8  = @(0.65843621399177,  0.79672131147541,  0.900473933649289, 0.786290322580645, 0.695966907962771)
  9  = @(0.844444444444444, 0.86084142394822,  0.800947867298578, 0.812,             0.813928182807399)
  10 = @(0.581395348837209, 0.533546325878594, 0.623853211009174, 0.690839694656489, 0.619895287958115)
  11 = @(0.79746835443038,  0.776073619631902, 0.829596412556054, 0.8,               0.665938864628821)
  12 = @(0.760504201680672, 0.778156996587031, 0.831683168316832, 0.858921161825726, 0.798553719008264)
  13 = @(0.839449541284404, 0.8561872909699,   0.854077253218884, 0.772549019607843, 0.717127071823204)
  14 = @(0.780392156862745, 0.853820598006645, 0.847290640394089, 0.792307692307692, 0.793926247288503)
  15 = @(0.788235294117647, 0.777003484320557, 0.84549356223176,  0.743295019157088, 0.725690890481064)
  16 = @(0.763565891472868, 0.752380952380952, 0.722943722943723, 0.493877551020408, 0.5917225950783)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $vals[$i]
    }
}
